$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81 (shifts existing rows 81..175 down to 82..176)
$ws.Rows(81).Insert()

# Populate the new row 81 with the new record's data
$ws.Cells.Item(81, 1).Value = 4
$ws.Cells.Item(81, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(81, 3).Value = "Los Lagos"
$ws.Cells.Item(81, 4).Value = 45225
$ws.Cells.Item(81, 5).Value = 10
$ws.Cells.Item(81, 6).Value = 100112026
$ws.Cells.Item(81, 7).Value = "Haba"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 40
$ws.Cells.Item(81, 11).Value = 16000
$ws.Cells.Item(81, 12).Value = 16000
$ws.Cells.Item(81, 13).Value = 16000
$ws.Cells.Item(81, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(81, 15).Value = "Región Metropolitana"
$ws.Cells.Item(81, 16).Value = 640
$ws.Cells.Item(81, 17).Value = 25
$ws.Cells.Item(81, 18).Value = "Hortaliza"
